# Reverted to "home button" commit: restores the smaller feature-checklist
# table (A1:H9) with updated model/feature names and checkbox layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 (header labels) ---
$ws.Range("C1").Value = "newModel"
$ws.Range("D1").Value = "newModel2"
$ws.Range("E1").Value = "sdfcsfs"
$ws.Range("F1").Value = "adasd"
$ws.Range("G1").Value = "model4"
# H1 stays "save" (unchanged); clear the now-unused I1:O1 tail.
$ws.Range("I1:O1").ClearContents()

# --- Row 2: Sex -> WBC ---
$ws.Range("A2").Value = "WBC"
$ws.Range("C2").Value = 1

# --- Row 3: Species -> LYMF ---
$ws.Range("A3").Value = "LYMF"
$ws.Range("C3").Value = 1

# --- Row 4: WBC -> RBC ---
$ws.Range("A4").Value = "RBC"
$ws.Range("D4").ClearContents()
$ws.Range("C4").Value = 1
$ws.Range("E4").Value = 1
$ws.Range("H4").Value = 1
$ws.Range("I4").ClearContents()
$ws.Range("K4").ClearContents()

# --- Row 5: LYMF -> HGB ---
$ws.Range("A5").Value = "HGB"
$ws.Range("D5").ClearContents()
$ws.Range("O5").ClearContents()
$ws.Range("C5").Value = 1
$ws.Range("G5").Value = 1
$ws.Range("H5").Value = 1

# --- Row 6: GRAN -> MCH ---
$ws.Range("A6").Value = "MCH"
$ws.Range("E6").ClearContents()
$ws.Range("F6").Value = 1
$ws.Range("J6").ClearContents()
$ws.Range("L6").ClearContents()
$ws.Range("H6").Value = 1

# --- Row 7: MID -> MCHC ---
$ws.Range("A7").Value = "MCHC"
$ws.Range("C7").ClearContents()
$ws.Range("N7").ClearContents()
$ws.Range("H7").Value = 1

# --- Row 8: RBC -> MPV ---
$ws.Range("A8").Value = "MPV"
$ws.Range("H8").ClearContents()
$ws.Range("M8").ClearContents()
$ws.Range("D8").Value = 1
$ws.Range("E8").Value = 1

# --- Row 9: HGB -> PLT ---
$ws.Range("A9").Value = "PLT"
$ws.Range("D9").Value = 1

# --- Rows 10-13 (MCH, MCHC, MPV, PLT) no longer exist; remove them ---
$ws.Rows("10:13").ClearContents()

# Restore the original ("home button") selection/view state.
$ws.Range("A1").Select() | Out-Null
